# Staff excel update: correct the salutations/prefixes for the three
# newly added Assistant Professors in the STAFF-DATA/004 sheet so their
# names are consistent with the rest of the table (e.g. "Dr. ..", "Mr. ..").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value  = "Mrs. FRIEDA F S"
$ws.Range("A10").Value = "Mrs. PAVEENA A"
$ws.Range("A11").Value = "Mr. SELVAPRASANTH P"

# Leave the cursor on the name column of the last edited row, matching
# the saved selection state of the workbook.
$ws.Range("A11").Select()
